# The title of slide 1 ("So-Show") had a typo / awkward hyphenation.
# Rename it to a single clean phrase: "The Soshow".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Title
$titleShape.TextFrame.TextRange.Text = "The Soshow"
